# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Note: several "Price" values look numeric (e.g. 379.46) but must stay as
# literal text to preserve the source's exact formatting (fixed decimals,
# leading zeros, etc.), so those are written with a leading apostrophe
# ('' inside a single-quoted PowerShell string == a literal ') which makes
# Excel store them as text instead of auto-converting to a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.171.00'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '2.957.90'
$ws.Range('E3').Value = '  +2.51%  '
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''379.46'
$ws.Range('E5').Value = '  +3.65%  '
$ws.Range('D6').Value = '''104.89'
$ws.Range('E6').Value = '  +3.30%  '
$ws.Range('D7').Value = '''0.540'
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('E8').Value = '  -3.31%  '
$ws.Range('D9').Value = '''0.595'
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('D10').Value = '''37.10'
$ws.Range('E10').Value = '  +1.62%  '
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('D12').Value = '''0.0840'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').Value = '''18.43'
$ws.Range('E13').Value = '  +1.02%  '
$ws.Range('D14').Value = '3.415.69'
$ws.Range('E14').Value = '  +2.31%  '
$ws.Range('D15').Value = '''7.50'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '2.956.52'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').Value = '''0.965'
$ws.Range('E17').Value = '  +4.08%  '
$ws.Range('D18').Value = '51.075.32'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').Value = '''3.33'
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('D20').Value = '''7.38'
$ws.Range('E20').Value = '  +2.99%  '
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('D22').Value = '0.0₃0960'
$ws.Range('E22').Value = '  +2.29%  '
$ws.Range('D23').Value = '''69.45'
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').Value = '''261.66'
$ws.Range('E24').Value = '  +1.60%  '
$ws.Range('E25').Value = '  +5.60%  '
$ws.Range('D26').Value = '''7.63'
$ws.Range('E26').Value = '  +10.16%  '
$ws.Range('D27').Value = '''7.30'
$ws.Range('E27').Value = '  +22.86%  '
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = '''25.86'
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.111'
$ws.Range('E31').Value = '  +8.84%  '
$ws.Range('D32').Value = '''9.84'
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').Value = '''34.72'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('E34').Value = '  -2.19%  '
$ws.Range('D35').Value = '''51.06'
$ws.Range('E35').Value = '  +0.69%  '
$ws.Range('D36').Value = '''0.0448'
$ws.Range('E36').Value = '  +8.39%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').Value = '''3.07'
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('D39').Value = '''17.27'
$ws.Range('E39').Value = '  +2.50%  '
$ws.Range('D40').Value = '''2.58'
$ws.Range('E40').Value = '  -1.64%  '
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('E42').Value = '  +3.50%  '
$ws.Range('D43').Value = '''122.55'
$ws.Range('E43').Value = '  +4.09%  '
$ws.Range('D44').Value = '''22.03'
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('D45').Value = '''0.289'
$ws.Range('E45').Value = '  +23.90%  '
$ws.Range('D46').Value = '''2.06'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').Value = '''2.38'
$ws.Range('E47').Value = '  +2.74%  '
$ws.Range('D48').Value = '2.035.91'
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('D49').Value = '''3.22'
$ws.Range('E49').Value = '  +2.50%  '
$ws.Range('D50').Value = '''0.0347'
$ws.Range('E50').Value = '  +12.72%  '
$ws.Range('D51').Value = '''1.29'
$ws.Range('E51').Value = '  +3.24%  '
